$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the index card path for row 2 (P2)
$ws.Range("P2").Value = "/uploads/index_cards/indexcard-1768207654861-25492876.png"

# New case records appended as rows 12-21
$data = @(
    @(36, "NPS-III-2025-065", "Invalid Date", "Rainer Antiola", "Justin Merina", "78 Del Pilar St., Manila", "Physical Injuries", "", "Invalid Date", "Atty. P. Navarro", "25-1062", "MTC Br. 4, Manila", "Invalid Date", "Information filed", "Community service", "Injuries-2025"),
    @(37, "NPS-III-2025-099", "Invalid Date", "Kershua", "Jeffy", "Sitio Mabuhay, Cebu City", "Illegal Possession of Firearms", "", "Invalid Date", "Atty. S. Morales", "25-1070", "RTC Br. 12, Cebu", "Invalid Date", "Information filed", "N/A", "Firearms-2025"),
    @(38, "NPS-III-2025-088", "Invalid Date", "Miguel Ohara", "Peter Parker", "210 Luna St., Makati City", "Cyber Libel", "", "Invalid Date", "Atty. E. Flores", "N/A", "N/A", "Invalid Date", "Pending", "N/A", "Cybercrime-2025"),
    @(39, "NPS-III-2025-067", "Invalid Date", "Michael jackson", "Justin Bieber", "56 Chino Roces Ave., Makati", "Bouncing Checks (BP 22)", "", "Invalid Date", "Atty. M. Gomez", "25-1081", "MTC Br. 6, Makati", "Invalid Date", "Information filed", "Fine ₱20,000", "BP22-2025"),
    @(40, "No: NPS-III-2025-030", "Invalid Date", "Obiwan Kenobi", "Luke  Skywalker", "Barangay Sta. Cruz, Manila", "Murder", "", "Invalid Date", "Atty. J. Cruz", "25-1090", "RTC Br. 51, Manila", "Invalid Date", "Information filed", "Death Penalty", "QTheft-2025"),
    @(41, "NPS-III-2025-047", "Invalid Date", "Arvin Antiola", "Harry Potter", "89 Bonifacio St., Caloocan City", "Malicious Mischief", "", "Invalid Date", "Atty. L. Ramirez", "N/A", "N/A", "Invalid Date", "Case dismissed", "N/A", "Mischief-2026"),
    @(42, "NPS-III-2025-001", "Invalid Date", "Jakerax", "Kian Arcayan", "123 Mabini St., Quezon City", "Theft", "Invalid Date", "Invalid Date", "Atty. M. Bullanday", "25-1046", "MTC Br. 2, QC", "Invalid Date", "Information filed", "Fine 1 million", "Theft-2025"),
    @(43, "NPS-III-2025-031", "Invalid Date", "US Government", "Arvin Birada", "Barangay San Isidro of United States", "Violation of RA 9165", "Invalid Date", "Invalid Date", "Atty. L. Anrcayan", "25-1050", "RTC Br. 88, QC", "Invalid Date", "Information filed", "Death Penalty", "Drugs-2025"),
    @(44, "NPS-III-2025-000", "Invalid Date", "KIan Perez", "Diddy", "45 Rizal Ave., Manila", "Rape", "Invalid Date", "Invalid Date", "Atty. M. Galamiton", "N/A", "N/A", "Invalid Date", "Case dismissed", "N/A", "Estafa-2025"),
    @(45, "NPS-III-2025-009", "Invalid Date", "RIck", "Coral", "Purok 3, Barangay Malinis, Laguna", "Grave Threats", "Invalid Date", "Invalid Date", "Atty. R. Bullanday", "N/A", "N/A", "Invalid Date", "Case dismissed", "N/A", "Threats-2025")
)

$startRow = 12
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    if ($rec[7] -ne "") {
        $ws.Cells.Item($row, 8).Value = $rec[7]
    }
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $ws.Cells.Item($row, 11).Value = $rec[10]
    $ws.Cells.Item($row, 12).Value = $rec[11]
    $ws.Cells.Item($row, 13).Value = $rec[12]
    $ws.Cells.Item($row, 14).Value = $rec[13]
    $ws.Cells.Item($row, 15).Value = $rec[14]
    $ws.Cells.Item($row, 16).Value = $rec[15]
}
